$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 202
$ws1.Range("F6").Value = 511
$ws1.Range("F10").Value = 6811
$ws1.Range("F13").Value = 3131
$ws1.Range("F14").Value = 206
$ws1.Range("F17").Value = 553
$ws1.Range("F18").Value = 10

# Sheet "全部类型" (fourth sheet) - same rows duplicated, update accordingly
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 202
$ws4.Range("F8").Value = 511
$ws4.Range("F13").Value = 6811
$ws4.Range("F17").Value = 3131
$ws4.Range("F18").Value = 206
$ws4.Range("F21").Value = 553
$ws4.Range("F22").Value = 10
